$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 15).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 67 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = "None"
}
